# Add a new "CBM (kidney)" worksheet (with its data table) to the workbook,
# mirroring the structure of the existing "CBM (retina)"/"CBM (muscle)" sheets,
# and make it the active/selected sheet (as the last tab).

$wb = $excel.ActiveWorkbook

# --- 1. Create the new worksheet at the end of the tab strip -----------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add($null, $lastSheet)
$ws.Name = "CBM (kidney)"

# --- 2. Header row -------------------------------------------------------------
$ws.Range("A1").Value = "Reference"
$ws.Range("B1").Value = "Average"
$ws.Range("C1").Value = "SD"
$ws.Range("D1").Value = "N"
$ws.Range("E1").Value = "SE"

# --- 3. Data rows ----------------------------------------------------------------
$ws.Range("A2").Value = "Carlson et al., 2003"
$ws.Range("B2").Value = 178.16
$ws.Range("C2").Value = 35.61
$ws.Range("D2").Value = 19
$ws.Range("E2").Formula = "=C2/SQRT(D2)"

$ws.Range("A3").Value = "Chang et al., 2012"
$ws.Range("B3").Value = 224.2
$ws.Range("D3").Value = 4
$ws.Range("E3").Value = 27.7

$ws.Range("A4").Value = "Gambaro et al., 1992"
$ws.Range("B4").Value = 235.57
$ws.Range("E4").Value = 1.05

$ws.Range("A5").Value = "Yagihashi, 1978 (6 month-old rats)"
$ws.Range("B5").Value = 160.5
$ws.Range("C5").Value = 3.8
$ws.Range("D5").Value = 4
$ws.Range("E5").Formula = "=C5/SQRT(D5)"

$ws.Range("A6").Value = "Yagihashi, 1978 (7 month-old rats)"
$ws.Range("B6").Value = 184.6
$ws.Range("C6").Value = 6.5
$ws.Range("D6").Value = 4
$ws.Range("E6").Formula = "=C6/SQRT(D6)"

$ws.Range("A7").Value = "Ireland et al., 1977"
$ws.Range("B7").Value = 305
$ws.Range("E7").Value = 10

# --- 4. Turn the range into a table, matching the other sheets' style ------------
$lo = $ws.ListObjects.Add(1, $ws.Range("A1:E7"), $null, 1)
$lo.Name = "Table58"
$lo.TableStyle = "TableStyleLight1"

# --- 5. Column widths (best-fit, like the sibling sheets) ------------------------
$ws.Columns.Item(1).ColumnWidth = 23.998697916666668
$ws.Columns.Item(5).ColumnWidth = 11.330729166666666

# --- 6. Update the selection on "CBM (retina)" to span the whole table, as in
#        the target (done before activating the new sheet, since selecting a
#        range switches the active sheet to it) -----------------------------------
$retina = $wb.Worksheets.Item("CBM (retina)")
$retina.Range("A1:E4").Select() | Out-Null

# --- 7. Make the new sheet the active / selected tab, with E2 selected -----------
$ws.Activate() | Out-Null
$ws.Range("E2").Select() | Out-Null
